# Apply updated values to the ACF/RandomForest result sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.13
$ws.Range("A3").Value = -21.40600000000003
$ws.Range("C5").Value = -14.2788
$ws.Range("E5").Value = 13.03509999999999
$ws.Range("E9").Value = 13.73560000000001
$ws.Range("E11").Value = 13.34939999999999
$ws.Range("A14").Value = -20.49609999999998
$ws.Range("A21").Value = -21.36200000000001
$ws.Range("E21").Value = 12.66690000000001
$ws.Range("A23").Value = -21.49470000000003
$ws.Range("A25").Value = -22.52250000000003

$wb.Save()
